# Update the crypto tracker sheet with the latest "variations" snapshot.
# (commit: "fichier crypto avec les variations")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantités détenues (col D) ---
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 81
$ws.Range("D7").Value = 106
$ws.Range("D8").Value = 738
$ws.Range("D9").Value = 93
$ws.Range("D10").Value = 25

# --- Cours actuel (col I) ---
$ws.Range("I4").Value = 47010.35
$ws.Range("I5").Value = 1657.69
$ws.Range("I6").Value = 1.42
$ws.Range("I7").Value = 2.39
$ws.Range("I8").Value = 357.8
$ws.Range("I9").Value = 0.8167
$ws.Range("I10").Value = 0.1569

# --- Variations % (cols L:T), rows 4-10 ---
$ws.Range("L4").Value = -0.059357373475232435
$ws.Range("M4").Value = -0.06460256583708386
$ws.Range("N4").Value = 0.03109184085750627
$ws.Range("O4").Value = 0.0880861136118746
$ws.Range("P4").Value = 0.4630739176408079
$ws.Range("Q4").Value = 0.5619458247977889
$ws.Range("R4").Value = 4.0676405613264635
$ws.Range("S4").Value = 6.494903400209579
$ws.Range("T4").Value = 0.9778238178769939

$ws.Range("L5").Value = -0.07068063154096606
$ws.Range("M5").Value = 0.05271816359088433
$ws.Range("N5").Value = 0.1820692371549498
$ws.Range("O5").Value = 0.13441808662043842
$ws.Range("P5").Value = 0.1556354732714378
$ws.Range("Q5").Value = 0.6968970465409504
$ws.Range("R5").Value = 4.560949774605982
$ws.Range("S5").Value = 10.457542135724795
$ws.Range("T5").Value = 1.7513136513197305

$ws.Range("L6").Value = -0.0761686279812713
$ws.Range("M6").Value = -0.10904819656484875
$ws.Range("N6").Value = 0.14094498334364944
$ws.Range("O6").Value = 0.1768631652164136
$ws.Range("P6").Value = 0.12686435740099752
$ws.Range("Q6").Value = 2.755892585492447
$ws.Range("R6").Value = 3.5572808643531446
$ws.Range("S6").Value = 8.090988244199822
$ws.Range("T6").Value = 3.8635278667987945

$ws.Range("L7").Value = -0.148899392483318
$ws.Range("M7").Value = -0.1489200744752624
$ws.Range("N7").Value = 0.008774391410314024
$ws.Range("O7").Value = 0.298133493150404
$ws.Range("P7").Value = -0.19081011101116957
$ws.Range("Q7").Value = 2.8707062955515528
$ws.Range("R7").Value = 4.068462522611122
$ws.Range("S7").Value = 0.0
$ws.Range("T7").Value = 3.5475511066256313

$ws.Range("L8").Value = 0.0
$ws.Range("M8").Value = 0.0
$ws.Range("N8").Value = 0.0
$ws.Range("O8").Value = 0.0
$ws.Range("P8").Value = 0.0
$ws.Range("Q8").Value = 0.0
$ws.Range("R8").Value = 0.0
$ws.Range("S8").Value = 0.0
$ws.Range("T8").Value = 0.0

$ws.Range("L9").Value = -0.07065042465085242
$ws.Range("M9").Value = -0.10565179444555767
$ws.Range("N9").Value = 0.013811156711854978
$ws.Range("O9").Value = -0.21316258091036075
$ws.Range("P9").Value = 0.08802629658443152
$ws.Range("Q9").Value = 1.59006304410709
$ws.Range("R9").Value = 2.5049558625332877
$ws.Range("S9").Value = 17.100889352745263
$ws.Range("T9").Value = 2.059585875105434

$ws.Range("L10").Value = -0.24889287096410356
$ws.Range("M10").Value = -0.16987429918132363
$ws.Range("N10").Value = -0.1355646710294664
$ws.Range("O10").Value = 0.16641709008449596
$ws.Range("P10").Value = 1.757659711821924
$ws.Range("Q10").Value = 1.5119924022236177
$ws.Range("R10").Value = 0.22570789904291685
$ws.Range("S10").Value = 2.4461718499814844
$ws.Range("T10").Value = 2.2671463756539505

# --- Horodatage de mise à jour (col V) : tout le monde pointe vers la
#     dernière heure de rafraîchissement ---
$ws.Range("V4").Value = "07/04/21 21:42"
$ws.Range("V5").Value = "07/04/21 21:42"
$ws.Range("V6").Value = "07/04/21 21:42"
$ws.Range("V7").Value = "07/04/21 21:42"
$ws.Range("V8").Value = "07/04/21 21:42"
$ws.Range("V9").Value = "07/04/21 21:42"
$ws.Range("V10").Value = "07/04/21 21:42"

# Reflect where the author last clicked before saving.
$ws.Range("L4").Select()
